# CRUD operation of Learning Journey on 29MAY2025
# Adds a new "MERN" Learning Journey record (row 6) to Sheet1, mirroring
# the existing "React" (row 2) entry's description/image and marking it
# "valid", then leaves the selection on D15 as in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "MERN"
$ws.Range("B6").Value = "This is an framework used in WebDevelopment"
$ws.Range("C6").Value = "../Utility/LJPIC.png"
$ws.Range("D6").Value = "valid"

$ws.Range("D15").Select() | Out-Null
